$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45190 = 2023-09-21) that
# was bumped by two days (serial 45192 = 2023-09-23) for every data row
# (rows 2 through 493).
$lastRow = 493
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45192
